# Regenerate the localization-status report: flip every "Ready for handoff"
# status cell over to "In Translation" (workflow moved from handoff into the
# translation stage) and re-fit the status columns to the new, shorter
# content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$fitWidth = 13.4101845877511

# --- Overview sheet: zh-cn (E) and de-de (F) status columns, rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Columns("E:F").AutoFit() | Out-Null
$wsOverview.Columns("E:F").ColumnWidth = $fitWidth

# --- Per-locale detail sheets: Status column (C), rows 2-4 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2:C4").Value = $newStatus
    $ws.Columns("C:C").AutoFit() | Out-Null
    $ws.Columns("C:C").ColumnWidth = $fitWidth
}
